$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "27.207.71"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "1.852.92"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("E4").Value = "  -0.31%  "
Set-TextValue "D5" "313.26"
$ws.Range("E5").Value = "  +0.56%  "
Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  +0.16%  "
Set-TextValue "D8" "0.3720"
$ws.Range("E8").Value = "  +0.34%  "
Set-TextValue "D9" "0.07275"
$ws.Range("E9").Value = "  -0.79%  "
Set-TextValue "D10" "0.8867"
$ws.Range("E10").Value = "  +1.51%  "
Set-TextValue "D11" "20.05"
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.953.19"
$ws.Range("E12").Value = "  +6.11%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D13" "0.07813"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("E14").Value = "  +0.86%  "
Set-TextValue "D15" "6.516"
$ws.Range("E15").Value = "  -0.43%  "
Set-TextValue "D16" "90.92"
$ws.Range("E16").Value = "  -0.26%  "
Set-TextValue "D17" "1.002"
$ws.Range("E17").Value = "  -0.47%  "
Set-TextValue "D18" "0.000008913"
$ws.Range("E18").Value = "  +0.80%  "
Set-TextValue "D19" "1.001"
$ws.Range("E19").Value = "  -0.43%  "
Set-TextValue "D20" "14.73"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "27.233.32"
$ws.Range("E21").Value = "  +0.86%  "
Set-TextValue "D22" "5.065"
$ws.Range("E22").Value = "  -0.58%  "
Set-TextValue "D23" "10.49"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "2.136.39"
$ws.Range("E24").Value = "  +2.34%  "
Set-TextValue "D25" "1.950"
$ws.Range("E25").Value = "  +5.48%  "
Set-TextValue "D26" "151.66"
$ws.Range("E26").Value = "  -1.09%  "
Set-TextValue "D27" "18.38"
$ws.Range("E27").Value = "  -0.25%  "
Set-TextValue "D28" "2.037"
$ws.Range("E28").Value = "  -0.01%  "
Set-TextValue "D29" "115.58"
$ws.Range("E29").Value = "  +0.10%  "
Set-TextValue "D30" "5.056"
$ws.Range("E30").Value = "  -1.36%  "
Set-TextValue "D31" "0.08803"
$ws.Range("E31").Value = "  -1.07%  "
Set-TextValue "D32" "3.168"
$ws.Range("E32").Value = "  +7.11%  "
Set-TextValue "D33" "0.7662"
$ws.Range("E33").Value = "  +5.35%  "
Set-TextValue "D34" "1.168"
$ws.Range("E34").Value = "  +3.38%  "
Set-TextValue "D35" "4.499"
$ws.Range("E35").Value = "  +1.50%  "
Set-TextValue "D36" "2.728"
$ws.Range("E36").Value = "  +10.64%  "
Set-TextValue "D37" "1.095"
$ws.Range("E37").Value = "  +2.44%  "
Set-TextValue "D38" "0.01941"
$ws.Range("E38").Value = "  -0.19%  "
Set-TextValue "D39" "0.05219"
$ws.Range("E39").Value = "  -0.03%  "
Set-TextValue "D40" "2.934"
$ws.Range("E40").Value = "  -0.30%  "
Set-TextValue "D41" "7.024"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  -0.81%  "
Set-TextValue "D43" "0.1626"
$ws.Range("E43").Value = "  +0.43%  "
Set-TextValue "D44" "8.412"
$ws.Range("E44").Value = "  +3.01%  "
Set-TextValue "D45" "0.4786"
$ws.Range("E45").Value = "  -0.88%  "
Set-TextValue "D46" "10.34"
$ws.Range("E46").Value = "  +1.68%  "
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("E48").Value = "  -0.58%  "
Set-TextValue "D49" "1.636"
$ws.Range("E49").Value = "  +0.30%  "
Set-TextValue "D50" "0.06199"
$ws.Range("E50").Value = "  +0.15%  "
Set-TextValue "D51" "65.32"
$ws.Range("E51").Value = "  +0.65%  "
